$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 3-8 with new data
$ws.Range("A3").Value = "2022-09-06 22:53:48"
$ws.Range("B3").Value = "MER_CT_Bagamoio_24"
$ws.Range("C3").Value = "MER C&T"

$ws.Range("A4").Value = "2022-09-06 22:53:49"
$ws.Range("B4").Value = "MER_CT_Bagamoio_24"
$ws.Range("C4").Value = "MER C&T"

$ws.Range("A5").Value = "2022-09-06 22:53:49"
$ws.Range("B5").Value = "MER_CT_Bagamoio_24"
$ws.Range("C5").Value = "MER C&T"
$ws.Range("D5").Value = "Buscar valores para cada indicador: DSD TX NEW"

$ws.Range("A6").Value = "2022-09-06 22:53:52"
$ws.Range("B6").Value = "MER_CT_Bagamoio_24"
$ws.Range("C6").Value = "MER C&T"
$ws.Range("D6").Value = "Buscar valores para cada indicador: DSD TX CURR"
$ws.Range("E6").Value = "ok"

$ws.Range("A7").Value = "2022-09-06 22:53:56"
$ws.Range("B7").Value = "MER_CT_Bagamoio_24"
$ws.Range("C7").Value = "MER C&T"
$ws.Range("D7").Value = "Buscar valores para cada indicador: DSD TX RTT"
$ws.Range("E7").Value = "ok"

$ws.Range("A8").Value = "2022-09-06 22:54:01"
$ws.Range("B8").Value = "MER_CT_Bagamoio_24"
$ws.Range("C8").Value = "MER C&T"
$ws.Range("D8").Value = "Buscar valores para cada indicador: DSD TX ML"
$ws.Range("E8").Value = "ok"

# Add new rows 9-12
$ws.Range("A9").Value = "2022-09-06 22:54:14"
$ws.Range("B9").Value = "MER_CT_Bagamoio_24"
$ws.Range("C9").Value = "MER C&T"
$ws.Range("D9").Value = "Buscar valores para cada indicador: DSD PMCT ART"
$ws.Range("E9").Value = "ok"

$ws.Range("A10").Value = "2022-09-06 22:54:15"
$ws.Range("B10").Value = "MER_CT_Bagamoio_24"
$ws.Range("C10").Value = "MER C&T"
$ws.Range("D10").Value = "Buscar valores para cada indicador: DSD TX PVLS"
$ws.Range("E10").Value = "ok"

$ws.Range("A11").Value = "2022-09-06 22:54:21"
$ws.Range("B11").Value = "MER_CT_Bagamoio_24"
$ws.Range("C11").Value = "MER C&T"
$ws.Range("D11").Value = "Buscar valores para cada indicador: DSD TX TB"
$ws.Range("E11").Value = "ok"

$ws.Range("A12").Value = "2022-09-06 22:54:25"
$ws.Range("B12").Value = "MER_CT_Bagamoio_24"
$ws.Range("C12").Value = "MER C&T"
$ws.Range("D12").Value = "Buscar valores para cada indicador: DSD TB ART"
$ws.Range("E12").Value = "ok"
